$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SQL Scripts Assignment")

# ---------------------------------------------------------------------------
# 1) New analysis-idea rows (17-26): fill column C (and B) first so the new
#    shared strings land before the "Hutania" author-name string, matching
#    the order the strings were actually typed in.
# ---------------------------------------------------------------------------

$ws.Range("B17").Value = "Not Queried"
$ws.Range("C17").Value = "Distinct count of projects by country code and histogram it out in tableau"

$ws.Range("B18").Value = "Not Queried"
$ws.Range("C18").Value = "Take top 100 projects based off 100 highest 'ProportionReached' ratios to see which categories/subcategories were most common"

$ws.Range("B19").Value = "Not Queried"
$ws.Range("C19").Value = "Regression analysis on 'ProportionReached' ratio versus other variables and see if there is an impact (variables: 'Backers' 'Goal Amount')"

$ws.Range("B20").Value = "Not Queried"
$ws.Range("C20").Value = "Group by year, then look at 'ProjectState' attribute to see breakdown of successful/failed/cancelled and if there is a noticeable difference as time progressed (visualize this in tableau)"

$ws.Range("B21").Value = "Not Queried"
$ws.Range("C21").Value = "Sucessful projects that had 100, or 50 or fewer backers. Compared to projects with more backers. This analysis would show a trend in (do larger groups lead to more successful projects, does category/subcategory have any impact on number of backers?)"

$ws.Range("B22").Value = "Not Queried"
$ws.Range("C22").Value = "Analysis on active days versus successful projects. Take min and take max, divide into two bins, compare successful/failed/cancelled of two bins"

$ws.Range("B23").Value = "Not Queried"
$ws.Range("C23").Value = "What day of the week averages the highest pledge amounts across all projects"

$ws.Range("B24").Value = "Not Queried"
$ws.Range("C24").Value = "Take top 10 projects based off whatever you please, divide each project into thirds based off time, and do analysis on what majority of money came in. Is there any trend. (visualize this in tableau)"

$ws.Range("B25").Value = "Not Queried"
$ws.Range("C25").Value = "Group by category, do majority of projects fall into specific subcategories etc."

$ws.Range("B26").Value = "Not Queried"
$ws.Range("C26").Value = "Group project by country, then visualize counts in tableau"

# Author column (A) filled in last, after all the descriptions -> new
# shared string "Hutania" ends up last in sharedStrings.xml, matching source.
$ws.Range("A17").Value = "Hutania"
$ws.Range("A18").Value = "Hutania"
$ws.Range("A19").Value = "Hutania"
$ws.Range("A20").Value = "Hutania"
$ws.Range("A21").Value = "Hutania"
$ws.Range("A22").Value = "Hutania"
$ws.Range("A23").Value = "Hutania"
$ws.Range("A24").Value = "Hutania"
$ws.Range("A25").Value = "Hutania"
$ws.Range("A26").Value = "Hutania"

# ---------------------------------------------------------------------------
# 2) Column C got "Wrap Text" turned on for the whole used range (header +
#    all data rows), which is why C1 (bold header) and C2:C26 pick up new
#    wrap-text cell styles while A/B stay on the plain default style.
# ---------------------------------------------------------------------------
$ws.Range("C1:C26").WrapText = $true

# New rows (17-26) were typed with the workbook's base Calibri 11 font
# instead of the sheet's customary size-12 font - only column C was touched
# so only column C picks up the new font size.
$ws.Range("C17:C26").Font.Size = 11

# ---------------------------------------------------------------------------
# 3) A couple of the longer descriptions needed their rows resized so the
#    wrapped text would display fully.
# ---------------------------------------------------------------------------
$ws.Rows.Item(4).RowHeight = 31
$ws.Rows.Item(14).RowHeight = 46.5
$ws.Rows.Item(18).RowHeight = 29
$ws.Rows.Item(19).RowHeight = 29

# ---------------------------------------------------------------------------
# 4) Selection moved to C28 as the last thing the author did.
# ---------------------------------------------------------------------------
$ws.Range("C28").Select() | Out-Null
